$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.862.54"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "3.166.44"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'568.19"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").Value = "'167.64"
$ws.Range("E6").Value = "  -6.77%  "
$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  -6.32%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.166.27"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "'0.120"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("D11").Value = "'6.69"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'0.384"
$ws.Range("E12").Value = "  -4.99%  "
$ws.Range("D13").Value = "3.726.76"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "63.933.01"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "'25.26"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").Value = "'0.0000158"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "3.166.88"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "'414.00"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "'5.33"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "'12.75"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'7.09"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'70.74"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "'0.203"
$ws.Range("E25").Value = "  +3.66%  "
$ws.Range("D26").Value = "'0.490"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").Value = "'0.0000106"
$ws.Range("E27").Value = "  -5.71%  "
$ws.Range("D28").Value = "'8.69"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'1.84"
$ws.Range("E30").Value = "  -6.41%  "
$ws.Range("D31").Value = "'21.75"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'4.98"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "'6.34"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("D35").Value = "'1.13"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("D36").Value = "'155.85"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").Value = "'1.37"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").Value = "2.737.68"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'1.70"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "'24.89"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("D41").Value = "'4.16"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("D42").Value = "'0.712"
$ws.Range("E42").Value = "  -7.22%  "
$ws.Range("D43").Value = "'38.64"
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0623"
$ws.Range("E44").Value = "  -5.27%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.64"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0262"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("D47").Value = "'296.17"
$ws.Range("E47").Value = "  -5.91%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'21.77"
$ws.Range("E48").Value = "  -5.68%  "
$ws.Range("D49").Value = "'2.05"
$ws.Range("E49").Value = "  -10.76%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0991"
$ws.Range("E51").Value = "  -6.32%  "
